$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-13 (columns D, J, K, L, M, P change; other columns
# stay the same for all rows, and row 13 is new with the same fixed columns).
$data = @(
    @{Row=2;  D=44377; J=650; K=14000; L=15000; M=14538; P=1454},
    @{Row=3;  D=44204; J=400; K=10000; L=11000; M=10500; P=1050},
    @{Row=4;  D=44160; J=360; K=10000; L=11000; M=10500; P=1050},
    @{Row=5;  D=44330; J=300; K=13000; L=14000; M=13500; P=1350},
    @{Row=6;  D=44406; J=400; K=14000; L=15000; M=14500; P=1450},
    @{Row=7;  D=44265; J=200; K=15000; L=16000; M=15500; P=1550},
    @{Row=8;  D=44358; J=300; K=14000; L=15000; M=14500; P=1450},
    @{Row=9;  D=44218; J=320; K=10000; L=11000; M=10500; P=1050},
    @{Row=10; D=44291; J=200; K=13000; L=14000; M=13500; P=1350},
    @{Row=11; D=44263; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=12; D=44441; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=13; D=44428; J=300; K=15000; L=16000; M=15500; P=1550}
)

foreach ($entry in $data) {
    $r = $entry.Row

    if ($r -eq 13) {
        # Row 13 is brand new; fill in all the fixed columns that are the
        # same for every record, then the variable ones below.
        $ws.Cells.Item($r, 1).Value = 1
        $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
        $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
        $ws.Cells.Item($r, 5).Value = 15
        $ws.Cells.Item($r, 6).Value = 100112003
        $ws.Cells.Item($r, 7).Value = "Ajo"
        $ws.Cells.Item($r, 8).Value = "Chino"
        $ws.Cells.Item($r, 9).Value = "Primera"
        $ws.Cells.Item($r, 14).Value = "`$/caja 10 kilos"
        $ws.Cells.Item($r, 15).Value = "China"
        $ws.Cells.Item($r, 17).Value = 10
        $ws.Cells.Item($r, 18).Value = "Hortaliza"

        # Copy the style of the row above for consistent formatting (D col
        # date style in particular).
        $ws.Range("A12:R12").Copy() | Out-Null
        $ws.Range("A13:R13").PasteSpecial(-4122) | Out-Null
    }

    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 16).Value = $entry.P
}
